$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 03:55:28"
$ws1.Range("A3").Value = "Total filas: 2"

# Update existing data row 6 (scrape time, arrival time, minutes changed;
# Linea/Parada columns unchanged)
$ws1.Range("A6").Value = "03:55:28"
$ws1.Range("B6").Value = "04:46"
$ws1.Range("D6").Value = 51

# New data row 7
$ws1.Range("A7").Value = "03:55:28"
$ws1.Range("B7").Value = "05:39"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 104
$ws1.Range("E7").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:55:28"

$ws2.Range("A6").Value = "03:55:28"
$ws2.Range("B6").Value = "04:46"
$ws2.Range("D6").Value = 51

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:55:28"
